$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "58.954.27"
$ws.Range("E2").Value = "  +2.09%  "
Set-TextValue $ws.Range("D3") "2.584.57"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "521.65"
$ws.Range("E5").Value = "  +1.15%  "
Set-TextValue $ws.Range("D6") "139.01"
$ws.Range("E6").Value = "  -2.04%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.01%  "
Set-TextValue $ws.Range("D9") "2.593.03"
$ws.Range("E9").Value = "  +0.46%  "
Set-TextValue $ws.Range("D10") "6.57"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("E13").Value = "  +3.37%  "
Set-TextValue $ws.Range("D14") "3.040.05"
$ws.Range("E14").Value = "  +0.75%  "
Set-TextValue $ws.Range("D15") "58.881.10"
$ws.Range("E15").Value = "  +2.07%  "
Set-TextValue $ws.Range("D16") "20.42"
$ws.Range("E16").Value = "  +1.10%  "
Set-TextValue $ws.Range("D17") "2.582.42"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("E18").Value = "  -0.43%  "
Set-TextValue $ws.Range("D19") "337.88"
$ws.Range("E19").Value = "  +0.36%  "
Set-TextValue $ws.Range("D20") "4.28"
$ws.Range("E20").Value = "  +0.14%  "
Set-TextValue $ws.Range("D21") "10.08"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("E23").Value = "  +0.03%  "
Set-TextValue $ws.Range("D24") "65.92"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").Value = "  +0.71%  "
Set-TextValue $ws.Range("D27") "0.998"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +0.64%  "
Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -2.10%  "
Set-TextValue $ws.Range("D31") "5.91"
$ws.Range("E31").Value = "  -4.59%  "
$ws.Range("E32").Value = "  +0.49%  "
Set-TextValue $ws.Range("D33") "18.67"
$ws.Range("E33").Value = "  +0.12%  "
Set-TextValue $ws.Range("D34") "149.04"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("E36").Value = "  -1.34%  "
Set-TextValue $ws.Range("D37") "36.78"
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("E38").Value = "  +1.94%  "
Set-TextValue $ws.Range("D39") "0.822"
$ws.Range("E39").Value = "  -0.84%  "
Set-TextValue $ws.Range("D40") "0.807"
$ws.Range("E40").Value = "  -7.27%  "
Set-TextValue $ws.Range("D41") "3.49"
Set-TextValue $ws.Range("D42") "0.997"
$ws.Range("E42").Value = "  +0.05%  "
Set-TextValue $ws.Range("D43") "271.19"
$ws.Range("E43").Value = "  +0.85%  "
Set-TextValue $ws.Range("D44") "10.74"
$ws.Range("E44").Value = "  +0.84%  "
Set-TextValue $ws.Range("D45") "0.0952"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("E47").Value = "  -0.69%  "
Set-TextValue $ws.Range("D48") "18.38"
$ws.Range("E48").Value = "  -1.69%  "
Set-TextValue $ws.Range("D49") "1.961.51"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  +0.06%  "
Set-TextValue $ws.Range("D51") "4.49"
$ws.Range("E51").Value = "  -2.11%  "
